$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value (applied as text to preserve exact formatting)
$updates = @(
    @{Cell="D2"; Value="26.293.20"}
    @{Cell="E2"; Value="  +0.62%  "}
    @{Cell="D3"; Value="1.601.34"}
    @{Cell="E3"; Value="  +1.14%  "}
    @{Cell="E4"; Value="  -0.01%  "}
    @{Cell="D5"; Value="212.50"}
    @{Cell="E5"; Value="  +0.63%  "}
    @{Cell="E6"; Value="  +0.30%  "}
    @{Cell="E7"; Value="  +0.04%  "}
    @{Cell="E8"; Value="  +0.31%  "}
    @{Cell="E9"; Value="  -0.19%  "}
    @{Cell="D10"; Value="18.96"}
    @{Cell="E10"; Value="  -1.12%  "}
    @{Cell="D11"; Value="0.0855"}
    @{Cell="E11"; Value="  +1.12%  "}
    @{Cell="D12"; Value="1.827.91"}
    @{Cell="E12"; Value="  +1.16%  "}
    @{Cell="D13"; Value="1.610.15"}
    @{Cell="E13"; Value="  +1.62%  "}
    @{Cell="E14"; Value="  +0.06%  "}
    @{Cell="D15"; Value="0.507"}
    @{Cell="E15"; Value="  -1.83%  "}
    @{Cell="D16"; Value="63.64"}
    @{Cell="E16"; Value="  -0.52%  "}
    @{Cell="D17"; Value="26.295.77"}
    @{Cell="E17"; Value="  +0.47%  "}
    @{Cell="D18"; Value="227.70"}
    @{Cell="E18"; Value="  +6.87%  "}
    @{Cell="D20"; Value="7.60"}
    @{Cell="E20"; Value="  +3.87%  "}
    @{Cell="E21"; Value="  -0.01%  "}
    @{Cell="E22"; Value="  +1.85%  "}
    @{Cell="E23"; Value="  +0.27%  "}
    @{Cell="D24"; Value="8.95"}
    @{Cell="E24"; Value="  +0.34%  "}
    @{Cell="D25"; Value="145.44"}
    @{Cell="E25"; Value="  +1.16%  "}
    @{Cell="E26"; Value="  +0.04%  "}
    @{Cell="E27"; Value="  -0.24%  "}
    @{Cell="E28"; Value="  +1.52%  "}
    @{Cell="E29"; Value="  +2.31%  "}
    @{Cell="E30"; Value="  -0.31%  "}
    @{Cell="E31"; Value="  +1.26%  "}
    @{Cell="E32"; Value="  +0.90%  "}
    @{Cell="D33"; Value="1.441.63"}
    @{Cell="E33"; Value="  +7.72%  "}
    @{Cell="E34"; Value="  +1.13%  "}
    @{Cell="E35"; Value="  -0.41%  "}
    @{Cell="D37"; Value="0.564"}
    @{Cell="E37"; Value="  -2.54%  "}
    @{Cell="E38"; Value="  -0.94%  "}
    @{Cell="E39"; Value="  +1.03%  "}
    @{Cell="D40"; Value="5.81"}
    @{Cell="E40"; Value="  +0.59%  "}
    @{Cell="E41"; Value="  +0.15%  "}
    @{Cell="E42"; Value="  +2.10%  "}
    @{Cell="D43"; Value="0.924"}
    @{Cell="E43"; Value="  -1.88%  "}
    @{Cell="D44"; Value="1.739.67"}
    @{Cell="E44"; Value="  +1.15%  "}
    @{Cell="D45"; Value="0.759"}
    @{Cell="E45"; Value="  -0.88%  "}
    @{Cell="E46"; Value="  -0.03%  "}
    @{Cell="D47"; Value="87.66"}
    @{Cell="E47"; Value="  +2.00%  "}
    @{Cell="E48"; Value="  +0.70%  "}
    @{Cell="D49"; Value="0.0500"}
    @{Cell="E49"; Value="  -0.14%  "}
    @{Cell="B50"; Value="Algorand"}
    @{Cell="C50"; Value="https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"}
    @{Cell="D50"; Value="0.0952"}
    @{Cell="E50"; Value="  -3.13%  "}
    @{Cell="B51"; Value="USDD"}
    @{Cell="C51"; Value="https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"}
    @{Cell="D51"; Value="0.999"}
    @{Cell="E51"; Value="  +0.04%  "}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
